$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 99.98 -> 0M
$t.Cell(1, 1).Range.Text = "0M"

# Row 2: 0.04 -> 0M
$t.Cell(2, 1).Range.Text = "0M"

# Row 3: 279 -> 0M
$t.Cell(3, 1).Range.Text = "0M"

# Insert 10 new rows after row 3 with the given values
$newValues = @('1023', '0.00003', '0.00012', '0.00004', '0.00001', '0.00004', '0.00004', '0.00005', '0.04325', '100.0')
$insertBefore = 4
foreach ($val in $newValues) {
    $t.Rows.Add($t.Rows.Item($insertBefore)) | Out-Null
    $t.Cell($insertBefore, 1).Range.Text = $val
    $insertBefore++
}

# The three multi-run rows (with tab-separated stats) collapse to single summary values.
# They are the last 3 rows of the table (originally rows 34/35/36, now shifted by +10).
$rowCount = $t.Rows.Count
$t.Cell($rowCount - 2, 1).Range.Text = "99.98"
$t.Cell($rowCount - 1, 1).Range.Text = "0.04"
$t.Cell($rowCount, 1).Range.Text = "279"
